$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "68.195.91"
$ws.Range("E2").Value = "  +2.09%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.142.50"
$ws.Range("E3").Value = "  +2.48%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.15%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "579.24"
$ws.Range("E5").Value = "  +0.71%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "180.73"
$ws.Range("E6").Value = "  +7.52%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.999"
$ws.Range("E7").Value = "  -0.05%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.139.39"
$ws.Range("E8").Value = "  +2.48%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.520"
$ws.Range("E9").Value = "  +1.77%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.54"
$ws.Range("E10").Value = "  +2.54%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.154"
$ws.Range("E11").Value = "  +2.76%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.473"
$ws.Range("E12").Value = "  +1.04%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000244"
$ws.Range("E13").Value = "  +2.10%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "37.04"
$ws.Range("E14").Value = "  +4.16%  "

$ws.Range("E15").Value = "  +0.92%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.657.68"
$ws.Range("E16").Value = "  +2.25%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "68.050.13"
$ws.Range("E17").Value = "  +1.99%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "7.12"
$ws.Range("E18").Value = "  +1.86%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.131.46"
$ws.Range("E19").Value = "  +2.27%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "16.56"
$ws.Range("E20").Value = "  -1.41%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "489.95"
$ws.Range("E21").Value = "  -0.19%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.699"
$ws.Range("E22").Value = "  +1.68%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.78"
$ws.Range("E23").Value = "  +1.65%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "84.06"
$ws.Range("E24").Value = "  +1.61%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.35"
$ws.Range("E25").Value = "  +7.30%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "12.91"
$ws.Range("E26").Value = "  +1.86%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.68"
$ws.Range("E27").Value = "  +5.43%  "

$ws.Range("E28").Value = "  +0.10%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.17"
$ws.Range("E29").Value = "  +5.45%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.37"
$ws.Range("E30").Value = "  +4.86%  "

$ws.Range("E31").Value = "  +1.14%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "28.40"
$ws.Range("E32").Value = "  +3.45%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.114"
$ws.Range("E33").Value = "  +2.01%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0₃0957"
$ws.Range("E34").Value = "  +5.44%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.998"
$ws.Range("E35").Value = "  -0.18%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "49.09"
$ws.Range("E36").Value = "  +5.50%  "

$ws.Range("B37").Value = "Filecoin"
$ws.Range("C37").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.69"
$ws.Range("E37").Value = "  +1.93%  "

$ws.Range("B38").Value = "Mantle"
$ws.Range("C38").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.959"
$ws.Range("E38").Value = "  +1.28%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.326"
$ws.Range("E39").Value = "  +8.85%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.06"
$ws.Range("E40").Value = "  +4.67%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "49.43"
$ws.Range("E41").Value = "  +0.57%  "

$ws.Range("E42").Value = "  +1.61%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "8.42"
$ws.Range("E43").Value = "  +1.34%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.72"
$ws.Range("E44").Value = "  +10.31%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "390.30"
$ws.Range("E45").Value = "  +6.65%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.788.97"
$ws.Range("E46").Value = "  +1.37%  "

$ws.Range("B47").Value = "InjectiveProtocol"
$ws.Range("C47").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "27.05"
$ws.Range("E47").Value = "  +11.05%  "

$ws.Range("B48").Value = "VeChain"
$ws.Range("C48").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0350"
$ws.Range("E48").Value = "  +1.44%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "136.75"
$ws.Range("E49").Value = "  +0.62%  "

$ws.Range("E51").Value = "  +8.48%  "
